$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column widths for the refreshed Nowcasts 2025Q4 table (B:K)
$ws.Columns.Item(2).ColumnWidth = 12.833333333333332
$ws.Columns.Item(3).ColumnWidth = 12.833333333333332
$ws.Columns.Item(4).ColumnWidth = 13.333333333333332
$ws.Columns.Item(5).ColumnWidth = 13.833333333333332
$ws.Columns.Item(6).ColumnWidth = 14.333333333333332
$ws.Columns.Item(7).ColumnWidth = 13.833333333333332
$ws.Columns.Item(8).ColumnWidth = 14.333333333333332
$ws.Columns.Item(9).ColumnWidth = 15.333333333333334
$ws.Columns.Item(10).ColumnWidth = 14.166666666666666
$ws.Columns.Item(11).ColumnWidth = 13.833333333333332
# Update header row (row 1) - same labels, now backed by the new Nowcasts 2025Q4 data block
$ws.Range("A1").Value = "Row"
$ws.Range("B1").Value = "Prognose"
$ws.Range("C1").Value = "surveys"
$ws.Range("D1").Value = "production"
$ws.Range("E1").Value = "orders"
$ws.Range("F1").Value = "turnover"
$ws.Range("G1").Value = "financial"
$ws.Range("H1").Value = "labor market"
$ws.Range("I1").Value = "prices"
$ws.Range("J1").Value = "national accounts"
$ws.Range("K1").Value = "Revision"

# Refresh rows 2-7 with the new 2025Q4 nowcast values (rows 8-11 are left untouched)
$ws.Range("B2").Value = 0.23819312193243655
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("A3").Value = "'2025-10-15"
$ws.Range("B3").Value = 0.27250757174013873
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = -0.04446289177310081
$ws.Range("E3").Value = 0.026866372944277718
$ws.Range("F3").Value = 0.013699165364939371
$ws.Range("G3").Value = 0.0077164346514254616
$ws.Range("H3").Value = 0.0017974243087641941
$ws.Range("I3").Value = -0.00035917559872720691
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0.029057119910123469
$ws.Range("A4").Value = "'2025-10-30"
$ws.Range("B4").Value = 0.47530117769821245
$ws.Range("C4").Value = 0.013819892787387423
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 0.0010814035564391603
$ws.Range("F4").Value = -0.0017173596431209128
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = -0.0031942376652086121
$ws.Range("I4").Value = 0.17924497381788901
$ws.Range("J4").Value = 0.011165551028892251
$ws.Range("K4").Value = 0.0023933820757954249
$ws.Range("A5").Value = "'2025-11-15"
$ws.Range("B5").Value = 0.45333859580121239
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 0.025413238481473435
$ws.Range("E5").Value = 0.0098231931439007394
$ws.Range("F5").Value = -0.077342023295886964
$ws.Range("G5").Value = -0.010602110093451271
$ws.Range("H5").Value = -0.0049827694136908461
$ws.Range("I5").Value = 0.017517360079840501
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 0.01821052920081434
$ws.Range("A6").Value = "'2025-11-30"
$ws.Range("B6").Value = 0.048661102019954994
$ws.Range("C6").Value = -0.33703266823001932
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 0.0025689436485147382
$ws.Range("F6").Value = -0.0022626448498755886
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = -0.0047813974867500367
$ws.Range("I6").Value = -0.064992962437577725
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 0.0018232355744505768
$ws.Range("A7").Value = "'2025-12-15"
$ws.Range("B7").Value = 0.096311646940341472
$ws.Range("C7").Value = 0
$ws.Range("D7").Value = 0.20861279869619584
$ws.Range("E7").Value = -0.014933617941889361
$ws.Range("F7").Value = -0.11620828574030728
$ws.Range("G7").Value = 0.0049454368109285657
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = -0.034765786904541268